$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The F column used to hold hyperlinks whose *display* text was a short
# "filename.png (WxH) (raw.githubusercontent.com)" caption; it now just
# holds the full raw.githubusercontent.com URL as plain text, and the
# hyperlink objects themselves are gone.
$ws.Hyperlinks.Delete()

$ws.Range("F1").Value = "Img"
$ws.Range("F2").Value = "https://raw.githubusercontent.com/hvijay31/Diet/main/brrice.png"
$ws.Range("F3").Value = "https://raw.githubusercontent.com/hvijay31/Diet/main/rice.png"
$ws.Range("F4").Value = "https://raw.githubusercontent.com/hvijay31/Diet/main/roti.png"
$ws.Range("F5").Value = "https://raw.githubusercontent.com/hvijay31/Diet/main/idili.png"
$ws.Range("F6").Value = "https://raw.githubusercontent.com/hvijay31/Diet/main/mutton.png"
$ws.Range("F7").Value = "https://raw.githubusercontent.com/hvijay31/Diet/main/chicken.png"
$ws.Range("F8").Value = "https://raw.githubusercontent.com/hvijay31/Diet/main/dosa.png"

# Widen column F (20 chars) now that it carries full URLs instead of short
# captions. 19 + 1/6 compensates for the fixed padding this host adds when
# converting the character-based ColumnWidth into the stored OOXML width,
# landing exactly on width="20".
$ws.Columns.Item(6).ColumnWidth = 19 + 1/6

Write-Host "Updated Img column links and widened column F"
